$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text so values like "1.000" or
# "0.05760" keep their exact formatting instead of being parsed as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.351.35"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.889.61"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.55"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4843"
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2902"
$ws.Range("E8").Value = "  -2.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06601"
$ws.Range("E9").Value = "  -2.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.896.18"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.92"
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07341"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.166"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.66"
$ws.Range("E14").Value = "  -1.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6618"
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.317.92"
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.42"
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007781"
$ws.Range("E18").Value = "  -2.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.448"
$ws.Range("E20").Value = "  +3.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.137.03"
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "193.39"
$ws.Range("E23").Value = "  -4.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.177"
$ws.Range("E24").Value = "  -1.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.385"
$ws.Range("E25").Value = "  -2.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.31"
$ws.Range("E26").Value = "  +1.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.19"
$ws.Range("E27").Value = "  -3.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.939"
$ws.Range("E28").Value = "  -1.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.451"
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.308"
$ws.Range("E30").Value = "  -0.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09131"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.044"
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05086"
$ws.Range("E33").Value = "  -9.04%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7311"
$ws.Range("E34").Value = "  -2.39%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.141"
$ws.Range("E35").Value = "  +1.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.706"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01792"
$ws.Range("E37").Value = "  -3.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.647"
$ws.Range("E38").Value = "  -2.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9198"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.070"
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.888"
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.67"
$ws.Range("E42").Value = "  -1.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4306"
$ws.Range("E43").Value = "  -4.40%  "
$ws.Range("E44").Value = "  -0.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.482"
$ws.Range("E45").Value = "  -3.28%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1328"
$ws.Range("E46").Value = "  -4.58%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.584"
$ws.Range("E47").Value = "  +9.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.93"
$ws.Range("E48").Value = "  -10.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.921"
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.95"
$ws.Range("E50").Value = "  -5.66%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05760"
$ws.Range("E51").Value = "  -4.24%  "
